# "Generate Report for Handoff":
# A new handoff was generated for the b10e2331-... source file, so its
# "Latest Handoff Datetime" is refreshed on both the zh-cn and de-de
# per-language status sheets.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 on each sheet corresponds to source file
# b10e2331-f6b8-4a8a-ab86-7f64f2d60f31.md ; column D is "Latest Handoff Datetime".
$wsZhCn.Range("D4").Value = "2016-03-07 01:59:00"
$wsDeDe.Range("D4").Value = "2016-03-07 01:59:10"
